$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the last data row (row 5), which held the "login cancel" /
# "already posted" messages. These shared strings become unused and are
# dropped from the workbook on save.
$ws.Range("A5:B5").ClearContents()

# Update the active selection to match the post-edit state.
$ws.Range("A5:B6").Select()
